$wb = $excel.ActiveWorkbook
$win = $excel.ActiveWindow
$win.Width = 14370
$win.Height = 6040
